$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A458:A463").EntireRow.Insert()

$ws.Cells.Item(458, 1).Value = 11
$ws.Cells.Item(458, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(458, 3).Value = "Bíobío"
$ws.Cells.Item(458, 4).Value = 45008
$ws.Cells.Item(458, 5).Value = 8
$ws.Cells.Item(458, 6).Value = "Fruta"
$ws.Cells.Item(458, 7).Value = 100103
$ws.Cells.Item(458, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(458, 9).Value = 100103006
$ws.Cells.Item(458, 10).Value = "Nectarín"
$ws.Cells.Item(458, 11).Value = "Artic Mist"
$ws.Cells.Item(458, 12).Value = "Especial"
$ws.Cells.Item(458, 13).Value = 250
$ws.Cells.Item(458, 14).Value = 17000
$ws.Cells.Item(458, 15).Value = 17000
$ws.Cells.Item(458, 16).Value = 17000
$ws.Cells.Item(458, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(458, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(458, 19).Value = 1062
$ws.Cells.Item(458, 20).Value = 16

$ws.Cells.Item(459, 1).Value = 11
$ws.Cells.Item(459, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(459, 3).Value = "Bíobío"
$ws.Cells.Item(459, 4).Value = 45008
$ws.Cells.Item(459, 5).Value = 8
$ws.Cells.Item(459, 6).Value = "Fruta"
$ws.Cells.Item(459, 7).Value = 100103
$ws.Cells.Item(459, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(459, 9).Value = 100103006
$ws.Cells.Item(459, 10).Value = "Nectarín"
$ws.Cells.Item(459, 11).Value = "Artic Mist"
$ws.Cells.Item(459, 12).Value = "Primera"
$ws.Cells.Item(459, 13).Value = 300
$ws.Cells.Item(459, 14).Value = 15000
$ws.Cells.Item(459, 15).Value = 15000
$ws.Cells.Item(459, 16).Value = 15000
$ws.Cells.Item(459, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(459, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(459, 19).Value = 938
$ws.Cells.Item(459, 20).Value = 16

$ws.Cells.Item(460, 1).Value = 11
$ws.Cells.Item(460, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(460, 3).Value = "Bíobío"
$ws.Cells.Item(460, 4).Value = 45008
$ws.Cells.Item(460, 5).Value = 8
$ws.Cells.Item(460, 6).Value = "Fruta"
$ws.Cells.Item(460, 7).Value = 100103
$ws.Cells.Item(460, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(460, 9).Value = 100103006
$ws.Cells.Item(460, 10).Value = "Nectarín"
$ws.Cells.Item(460, 11).Value = "Artic Mist"
$ws.Cells.Item(460, 12).Value = "Segunda"
$ws.Cells.Item(460, 13).Value = 250
$ws.Cells.Item(460, 14).Value = 12000
$ws.Cells.Item(460, 15).Value = 12000
$ws.Cells.Item(460, 16).Value = 12000
$ws.Cells.Item(460, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(460, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(460, 19).Value = 750
$ws.Cells.Item(460, 20).Value = 16

$ws.Cells.Item(461, 1).Value = 11
$ws.Cells.Item(461, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(461, 3).Value = "Bíobío"
$ws.Cells.Item(461, 4).Value = 45008
$ws.Cells.Item(461, 5).Value = 8
$ws.Cells.Item(461, 6).Value = "Fruta"
$ws.Cells.Item(461, 7).Value = 100103
$ws.Cells.Item(461, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(461, 9).Value = 100103006
$ws.Cells.Item(461, 10).Value = "Nectarín"
$ws.Cells.Item(461, 11).Value = "June Pearl"
$ws.Cells.Item(461, 12).Value = "Especial"
$ws.Cells.Item(461, 13).Value = 250
$ws.Cells.Item(461, 14).Value = 18000
$ws.Cells.Item(461, 15).Value = 18000
$ws.Cells.Item(461, 16).Value = 18000
$ws.Cells.Item(461, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(461, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(461, 19).Value = 1125
$ws.Cells.Item(461, 20).Value = 16

$ws.Cells.Item(462, 1).Value = 11
$ws.Cells.Item(462, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(462, 3).Value = "Bíobío"
$ws.Cells.Item(462, 4).Value = 45008
$ws.Cells.Item(462, 5).Value = 8
$ws.Cells.Item(462, 6).Value = "Fruta"
$ws.Cells.Item(462, 7).Value = 100103
$ws.Cells.Item(462, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(462, 9).Value = 100103006
$ws.Cells.Item(462, 10).Value = "Nectarín"
$ws.Cells.Item(462, 11).Value = "June Pearl"
$ws.Cells.Item(462, 12).Value = "Primera"
$ws.Cells.Item(462, 13).Value = 250
$ws.Cells.Item(462, 14).Value = 16000
$ws.Cells.Item(462, 15).Value = 16000
$ws.Cells.Item(462, 16).Value = 16000
$ws.Cells.Item(462, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(462, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(462, 19).Value = 1000
$ws.Cells.Item(462, 20).Value = 16

$ws.Cells.Item(463, 1).Value = 11
$ws.Cells.Item(463, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(463, 3).Value = "Bíobío"
$ws.Cells.Item(463, 4).Value = 45008
$ws.Cells.Item(463, 5).Value = 8
$ws.Cells.Item(463, 6).Value = "Fruta"
$ws.Cells.Item(463, 7).Value = 100103
$ws.Cells.Item(463, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(463, 9).Value = 100103006
$ws.Cells.Item(463, 10).Value = "Nectarín"
$ws.Cells.Item(463, 11).Value = "June Pearl"
$ws.Cells.Item(463, 12).Value = "Segunda"
$ws.Cells.Item(463, 13).Value = 200
$ws.Cells.Item(463, 14).Value = 14000
$ws.Cells.Item(463, 15).Value = 14000
$ws.Cells.Item(463, 16).Value = 14000
$ws.Cells.Item(463, 17).Value = "$/caja 16 kilos empedrada"
$ws.Cells.Item(463, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(463, 19).Value = 875
$ws.Cells.Item(463, 20).Value = 16
